# Updates cryptos list values (price & volume/1h columns) to the new snapshot.
# Leading apostrophe forces text entry (avoids numeric auto-conversion for
# values like "0.998" or "3.41"); Style reset keeps the cell on the default
# style index so no spurious style/number-format diff is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.607.48"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +4.33%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.371.18"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.70%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.23%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'593.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +6.47%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'186.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.20%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.08%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.598"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +3.89%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +5.31%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.586"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +1.57%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'47.30"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +3.30%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +7.35%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'641.31"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +12.38%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'3.904.88"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +1.41%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'8.53"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.16%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'68.639.65"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +4.30%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  +1.92%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.371.46"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.29%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'17.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +1.70%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  +2.64%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.911"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +2.53%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'17.97"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.32%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'5.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +1.79%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'99.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +1.80%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'4.09"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +3.99%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +6.28%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'9.77"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +4.55%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'32.91"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +8.14%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'8.69"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +2.93%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'6.84"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +1.89%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'612.69"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +9.54%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'3.73"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +1.16%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'3.979.97"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +6.48%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'11.10"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +2.65%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  +2.79%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -0.13%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'56.22"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +1.32%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'2.79"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +8.22%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +7.30%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.130"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +3.43%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'33.68"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.20%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  +3.18%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "'ApeXProtocol"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'3.41"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +1.86%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = "'TheGraph"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'0.344"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +3.43%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.0423"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +3.95%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  +2.86%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'2.59"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +3.37%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +0.40%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  +9.68%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'131.69"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +5.91%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'7.79"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +7.22%  "
$ws.Range("E51").Style = "Normal"
